$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits on the
#    "No because ..." paragraph. It will be re-created later at the end of
#    the rewritten "Problem 1" answer paragraph.
# ---------------------------------------------------------------------------
$d.Bookmarks.ShowHidden = $true
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Remove the paragraphs that get folded into the big rewritten answer:
#    old paragraph 5  -> "where I would use id as a primary, ..."
#    old paragraph 6  -> (empty paragraph)
#    old paragraph 7  -> "Describe using the following terminologies:"
#    old paragraphs 8-11 -> the four bulleted list items
#    These paragraphs sit right after paragraph 4 ("One called "Toy" ...")
#    and right before the "Problem 2" heading paragraph.
# ---------------------------------------------------------------------------
$pStart = $d.Paragraphs.Item(5)
$pEnd = $d.Paragraphs.Item(12)
$killRange = $d.Range($pStart.Range.Start, $pEnd.Range.Start)
$killRange.Delete()

# ---------------------------------------------------------------------------
# 3) Rewrite paragraph 3 ("I would model a one to many relationship ...")
#    with the new opening paragraph text.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$r3 = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$p1xml = '<w:r><w:t xml:space="preserve">I would store the list of toys in the database using </w:t></w:r><w:r><w:t>three table</w:t></w:r><w:r><w:t xml:space="preserve"> tables. One called Toy, which would have the fields of age, name and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>toy_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> as the primary key, and the other called Attributes, which would have </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>a</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>attribute_id</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> as primary ke</w:t></w:r><w:r><w:t>y</w:t></w:r><w:r><w:t xml:space="preserve"> and attributes as this tables field.</w:t></w:r><w:r><w:t xml:space="preserve"> The third table would be called </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>ToyAttributes</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">, which the primary key would be the mapping of one toy id as one Attributes id. </w:t></w:r><w:r><w:t xml:space="preserve">I would then map it so it would constitute of a many to </w:t></w:r><w:r><w:t>many</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>relationship</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> using the third table as the connection</w:t></w:r><w:r><w:t xml:space="preserve"> between the </w:t></w:r><w:r><w:t xml:space="preserve">other </w:t></w:r><w:r><w:t xml:space="preserve">two tables. </w:t></w:r>'
$r3.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $p1xml + '</w:p>')

# ---------------------------------------------------------------------------
# 4) Rewrite paragraph 4 ("One called "Toy" ...") with the new closing
#    paragraph text, and re-attach the "_GoBack" bookmark at its very end
#    (embedded directly in the inserted markup -- constructing a fresh
#    Range/Bookmarks.Add call after a structural edit is unreliable here).
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$r4 = $d.Range($p4.Range.Start, $p4.Range.End - 1)
$p2xml = '<w:r><w:t xml:space="preserve">This is so the tables are easy to read and maintainable as they are decoupled through the main information. The performance would be quick as you would be able to get all the attributes from the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>toy_attribute</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> table through the toy table (unless there is hundreds of records). The storage is quite efficient as the different attributes and toys are separated out, so they </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>is</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> the potential of sharing types of attributes between the toys, which minimise redundancy</w:t></w:r><w:r><w:t xml:space="preserve"> in the database</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">If new attributes are </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>added,  a</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> new way of mapping would have to be created between the toy associated and the new attribute. Overall, this way is the most maintainable, performance efficient, storage efficient and dynamic for insertion of new tables.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$r4.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $p2xml + '</w:p>')
